# Append new order-line rows (6-16) to the Webstaurant Bakery order sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "711MAVYELLOW"
$ws.Range("B6").Value = "Sprinkles - Yellow"
$ws.Range("C6").Value = "'1"
$ws.Range("D6").Value = "'59.99"
$ws.Range("E6").Value = "'59.99"

$ws.Range("A7").Value = "'10207579"
$ws.Range("B7").Value = "Salt - Sea Coarse"
$ws.Range("C7").Value = "'4"
$ws.Range("D7").Value = "'17.28"
$ws.Range("E7").Value = "'69.12"

$ws.Range("A8").Value = "245CBB25183"
$ws.Range("B8").Value = "Box Cake - Full Sheet (Bottom)"
$ws.Range("C8").Value = "'1"
$ws.Range("D8").Value = "'45.99"
$ws.Range("E8").Value = "'45.99"

$ws.Range("A9").Value = "409ML90266"
$ws.Range("B9").Value = "Choc Curls - Dark"
$ws.Range("C9").Value = "'3"
$ws.Range("D9").Value = "'113.35"
$ws.Range("E9").Value = "'340.05"

$ws.Range("A10").Value = "433SLINERBL"
$ws.Range("B10").Value = "Sheet Pan Liner - Silicone Coated"
$ws.Range("C10").Value = "'4"
$ws.Range("D10").Value = "'76.99"
$ws.Range("E10").Value = "'307.96"

$ws.Range("A11").Value = "5000TOUT96"
$ws.Range("B11").Value = "Java Box (96oz)"
$ws.Range("C11").Value = "'4"
$ws.Range("D11").Value = "'81.99"
$ws.Range("E11").Value = "'327.96"

$ws.Range("A12").Value = "245CCGR2518"
$ws.Range("B12").Value = "Cake Board - Full Sheet (SO)"
$ws.Range("C12").Value = "'1"
$ws.Range("D12").Value = "'56.99"
$ws.Range("E12").Value = "'56.99"

$ws.Range("A13").Value = "245CCGR1410BL"
$ws.Range("B13").Value = "Cake Board - 1/4 Sheet"
$ws.Range("C13").Value = "'1"
$ws.Range("D13").Value = "'32.70"
$ws.Range("E13").Value = "'32.70"

$ws.Range("A14").Value = "'150300865"
$ws.Range("B14").Value = "Bag Paper - 6x13.5 Window"
$ws.Range("C14").Value = "'4"
$ws.Range("D14").Value = "'79.99"
$ws.Range("E14").Value = "'319.96"

$ws.Range("A15").Value = "271241CUTC"
$ws.Range("B15").Value = "Mop Head Cut (White)"
$ws.Range("C15").Value = "'24"
$ws.Range("D15").Value = "'2.19"
$ws.Range("E15").Value = "'52.56"

$ws.Range("A16").Value = "5004CAFE"
$ws.Range("B16").Value = "Cup - Espresso (4oz)"
$ws.Range("C16").Value = "'1"
$ws.Range("D16").Value = "'34.99"
$ws.Range("E16").Value = "'34.99"

